$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with columns P (14) and Q (15), copying header style from O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update recalculated statistics in columns B:H for rows 2-25
$ws.Range("B2").Value = 3.592193807253011
$ws.Range("C2").Value = 0.9728239593682986
$ws.Range("D2").Value = 0.3515597420837508
$ws.Range("E2").Value = 1.413454056897237
$ws.Range("F2").Value = 6.643519541856051
$ws.Range("G2").Value = 0.0007750257785415873
$ws.Range("H2").Value = 0.01799275965849567

$ws.Range("B3").Value = 3.094737136111462
$ws.Range("C3").Value = 0.8385189313133594
$ws.Range("D3").Value = 0.3084519180545016
$ws.Range("E3").Value = 1.209056049506927
$ws.Range("F3").Value = 5.797049498982659
$ws.Range("G3").Value = 0.0007848095910370331
$ws.Range("H3").Value = 0.01098235956438964

$ws.Range("B4").Value = 2.794447868647239
$ws.Range("C4").Value = 0.7573603672388458
$ws.Range("D4").Value = 0.2826156348752562
$ws.Range("E4").Value = 1.086497025992742
$ws.Range("F4").Value = 5.289030495239814
$ws.Range("G4").Value = 0.0007909192209934907
$ws.Range("H4").Value = 0.007476543439597694

$ws.Range("B5").Value = 2.673157674641629
$ws.Range("C5").Value = 0.7245557031730812
$ws.Range("D5").Value = 0.2722186794152464
$ws.Range("E5").Value = 1.037177375852764
$ws.Range("F5").Value = 5.084493103539671
$ws.Range("G5").Value = 0.0007934375920689807
$ws.Range("H5").Value = 0.006226574508866101

$ws.Range("B6").Value = 2.653076913637335
$ws.Range("C6").Value = 0.7191230762707619
$ws.Range("D6").Value = 0.2704994688046298
$ws.Range("E6").Value = 1.029022424585079
$ws.Range("F6").Value = 5.050667260822962
$ws.Range("G6").Value = 0.0007938575755114414
$ws.Range("H6").Value = 0.00602920484120073

$ws.Range("B7").Value = 2.792808000874118
$ws.Range("C7").Value = 0.7569169410717791
$ws.Range("D7").Value = 0.2824749195000749
$ws.Range("E7").Value = 1.085829501060758
$ws.Range("F7").Value = 5.286262549375977
$ws.Range("G7").Value = 0.0007909530650476978
$ws.Range("H7").Value = 0.007458990086647188

$ws.Range("B8").Value = 3.419464390638552
$ws.Range("C8").Value = 0.9262062584412831
$ws.Range("D8").Value = 0.3365475539058593
$ws.Range("E8").Value = 1.342293390315518
$ws.Range("F8").Value = 6.348931269709993
$ws.Range("G8").Value = 0.0007783798507598116
$ws.Range("H8").Value = 0.01539576850346602

$ws.Range("B9").Value = 4.700372571961736
$ws.Range("C9").Value = 1.271674211786546
$ws.Range("D9").Value = 0.4490088317559042
$ws.Range("E9").Value = 1.874524697580242
$ws.Range("F9").Value = 8.54954363176509
$ws.Range("G9").Value = 0.0007543859231021604
$ws.Range("H9").Value = 0.03843096486116337

$ws.Range("B10").Value = 5.681851589069709
$ws.Range("C10").Value = 1.534933543203579
$ws.Range("D10").Value = 0.5145970482402049
$ws.Range("E10").Value = 2.181156984318235
$ws.Range("F10").Value = 10.01261651337683
$ws.Range("G10").Value = 0.0007377085274875015
$ws.Range("H10").Value = 0.05986773334596052

$ws.Range("B11").Value = 6.060852215372392
$ws.Range("C11").Value = 1.625702091616631
$ws.Range("D11").Value = 0.3624039816228475
$ws.Range("E11").Value = 1.458407103050831
$ws.Range("F11").Value = 8.520065099429274
$ws.Range("G11").Value = 0.0007365039337319664
$ws.Range("H11").Value = 0.07345247542934885

$ws.Range("B12").Value = 6.176123717061955
$ws.Range("C12").Value = 1.648966458223526
$ws.Range("D12").Value = 0.2453812846150498
$ws.Range("E12").Value = 0.9249065467925845
$ws.Range("F12").Value = 7.091009358015953
$ws.Range("G12").Value = 0.0007383314655146549
$ws.Range("H12").Value = 0.1058802324509287

$ws.Range("B13").Value = 6.095938445200261
$ws.Range("C13").Value = 1.622067523157682
$ws.Range("D13").Value = 0.149689935231109
$ws.Range("E13").Value = 0.5102820872318716
$ws.Range("F13").Value = 5.64520361269075
$ws.Range("G13").Value = 0.000742548817204121
$ws.Range("H13").Value = 0.1542586976334093

$ws.Range("B14").Value = 5.954518088575412
$ws.Range("C14").Value = 1.581582639859732
$ws.Range("D14").Value = 0.09610097634316617
$ws.Range("E14").Value = 0.2949742734336169
$ws.Range("F14").Value = 4.656702487624585
$ws.Range("G14").Value = 0.0007464675467949619
$ws.Range("H14").Value = 0.1983802917773119

$ws.Range("B15").Value = 5.878355865463732
$ws.Range("C15").Value = 1.560849980898581
$ws.Range("D15").Value = 0.0845554914513329
$ws.Range("E15").Value = 0.2506742056944766
$ws.Range("F15").Value = 4.39810619960187
$ws.Range("G15").Value = 0.0007480092756644297
$ws.Range("H15").Value = 0.2091117149792723

$ws.Range("B16").Value = 5.476528216380643
$ws.Range("C16").Value = 1.455231859980813
$ws.Range("D16").Value = 0.08257916679567501
$ws.Range("E16").Value = 0.23218824738451
$ws.Range("F16").Value = 4.141679871048069
$ws.Range("G16").Value = 0.0007540487573626753
$ws.Range("H16").Value = 0.1901492739375641

$ws.Range("B17").Value = 5.243095403503958
$ws.Range("C17").Value = 1.395375546797482
$ws.Range("D17").Value = 0.1065278786778023
$ws.Range("E17").Value = 0.3159033873218107
$ws.Range("F17").Value = 4.446211343870061
$ws.Range("G17").Value = 0.000756690802023105
$ws.Range("H17").Value = 0.1501246063140087

$ws.Range("B18").Value = 5.126723217684912
$ws.Range("C18").Value = 1.368304605023354
$ws.Range("D18").Value = 0.1667518495850544
$ws.Range("E18").Value = 0.5620078106294812
$ws.Range("F18").Value = 5.330518548724257
$ws.Range("G18").Value = 0.0007564029556796089
$ws.Range("H18").Value = 0.09872645696575688

$ws.Range("B19").Value = 5.117342187922588
$ws.Range("C19").Value = 1.371728868720254
$ws.Range("D19").Value = 0.2703118587963331
$ws.Range("E19").Value = 1.02543288285132
$ws.Range("F19").Value = 6.734736769748025
$ws.Range("G19").Value = 0.0007532705843329435
$ws.Range("H19").Value = 0.05929253359585118

$ws.Range("B20").Value = 5.419749786727493
$ws.Range("C20").Value = 1.464585242450767
$ws.Range("D20").Value = 0.4960211483651449
$ws.Range("E20").Value = 2.094196196223379
$ws.Range("F20").Value = 9.609822284096936
$ws.Range("G20").Value = 0.0007421174685810888
$ws.Range("H20").Value = 0.05368000980545684

$ws.Range("B21").Value = 6.209845449863394
$ws.Range("C21").Value = 1.678709832720926
$ws.Range("D21").Value = 0.5849397475758167
$ws.Range("E21").Value = 2.514291418895596
$ws.Range("F21").Value = 11.18677187458462
$ws.Range("G21").Value = 0.0007280661428483715
$ws.Range("H21").Value = 0.07568729869335122

$ws.Range("B22").Value = 6.738257924722575
$ws.Range("C22").Value = 1.821332885957645
$ws.Range("D22").Value = 0.6333957058354827
$ws.Range("E22").Value = 2.741253936927322
$ws.Range("F22").Value = 12.12027205584155
$ws.Range("G22").Value = 0.0007192112577410468
$ws.Range("H22").Value = 0.09106913576158249

$ws.Range("B23").Value = 6.454518149203068
$ws.Range("C23").Value = 1.74473614813877
$ws.Range("D23").Value = 0.6073203822525954
$ws.Range("E23").Value = 2.619194450518933
$ws.Range("F23").Value = 11.61837536592122
$ws.Range("G23").Value = 0.0007239461921763119
$ws.Range("H23").Value = 0.08266464601225998

$ws.Range("B24").Value = 5.414200745083178
$ws.Range("C24").Value = 1.464117913225209
$ws.Range("D24").Value = 0.5128304563538677
$ws.Range("E24").Value = 2.17544961145046
$ws.Range("F24").Value = 9.791084174102934
$ws.Range("G24").Value = 0.0007417274527224171
$ws.Range("H24").Value = 0.05477513114257526

$ws.Range("B25").Value = 4.346564607016091
$ws.Range("C25").Value = 1.176287933780657
$ws.Range("D25").Value = 0.4176814302950334
$ws.Range("E25").Value = 1.726490941428352
$ws.Range("F25").Value = 7.938087217534076
$ws.Range("G25").Value = 0.0007608289724725436
$ws.Range("H25").Value = 0.03122623659395796

# Populate the two new trailing columns (P, Q) with zeros for data rows
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}

